$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @(
    4367.620106690559,
    4293.968268859719,
    4168.197981837488,
    4168.197981837488,
    4168.197981837488,
    4128.671102361796,
    4128.671102361796,
    4105.729879787218,
    4105.729879787218,
    4105.729879787218,
    4105.729879787218
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
